$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 545; existing rows 545:560 shift down to 546:561
$ws.Rows.Item(545).Insert()

# Populate the newly inserted row 545 with the new weekly price record
$ws.Cells.Item(545, 1).Value = 6
$ws.Cells.Item(545, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(545, 3).Value = "Metropolitana"
$ws.Cells.Item(545, 4).Value = 44448
$ws.Cells.Item(545, 5).Value = 13
$ws.Cells.Item(545, 6).Value = 100112028
$ws.Cells.Item(545, 7).Value = "Sandia"
$ws.Cells.Item(545, 8).Value = "Sin especificar"
$ws.Cells.Item(545, 9).Value = "Primera"
$ws.Cells.Item(545, 10).Value = 4160
$ws.Cells.Item(545, 11).Value = 1200
$ws.Cells.Item(545, 12).Value = 1200
$ws.Cells.Item(545, 13).Value = 1200
$ws.Cells.Item(545, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(545, 15).Value = "Perú"
$ws.Cells.Item(545, 16).Value = 1200
$ws.Cells.Item(545, 17).Value = 1
$ws.Cells.Item(545, 18).Value = "Hortaliza"
